# Applies the "Updated symbol list" edit: refreshed prices/volumes for several
# coins and fixed a swapped CEJI / BKEXToken row pair (rows 42-43).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Single apostrophe used to force text-entry (like typing '271.81 into Excel)
# so the numeric-looking Price values stay text cells, matching the source data.
$q = [string][char]39

$ws.Range("D2").Value = $q + "271.81"
$ws.Range("D3").Value = $q + "22.78"
$ws.Range("D4").Value = $q + "6.339"
$ws.Range("D5").Value = $q + "0.06213"
$ws.Range("D6").Value = $q + "3.649"
$ws.Range("D7").Value = $q + "6.697"
$ws.Range("D8").Value = $q + "1.388"
$ws.Range("D9").Value = $q + "0.8298"
$ws.Range("D10").Value = $q + "0.01377"
$ws.Range("D11").Value = $q + "0.1608"
$ws.Range("D12").Value = $q + "0.08290"
$ws.Range("D13").Value = $q + "0.03442"
$ws.Range("D14").Value = $q + "0.03175"
$ws.Range("D15").Value = $q + "0.09342"
$ws.Range("D17").Value = $q + "0.001642"
$ws.Range("D18").Value = $q + "0.04715"
$ws.Range("D19").Value = $q + "0.006326"
$ws.Range("D20").Value = $q + "0.005674"
$ws.Range("D22").Value = $q + "0.0001499"
$ws.Range("D23").Value = $q + "3.719"
$ws.Range("D24").Value = $q + "2.399"
$ws.Range("D40").Value = $q + "0.04689"
$ws.Range("D41").Value = $q + "0.007020"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = $q + "0.1162"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = $q + "0.003348"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = $q + "0.01167"
$ws.Range("D45").Value = $q + "0.00006256"
$ws.Range("D46").Value = $q + "0.0009898"
$ws.Range("D48").Value = $q + "0.9199"
$ws.Range("D49").Value = $q + "0.002145"
$ws.Range("E50").Value = "49CryptobidCoinCBCWorstin24h"

# Drop the quote-prefix flag the forced-text entry above leaves behind, so the
# cell style matches a plain, never-formatted cell.
foreach ($addr in @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49")) {
    $ws.Range($addr).Style = "Normal"
}
